$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.200.30"
$ws.Range("E2").Value = "  -3.69%  "
$ws.Range("D3").Value = "2.974.44"
$ws.Range("E3").Value = "  -3.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "545.11"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.58"
$ws.Range("E6").Value = "  -7.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "2.974.15"
$ws.Range("E8").Value = "  -3.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.488"
$ws.Range("E9").Value = "  -2.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.97"
$ws.Range("E10").Value = "  -6.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.142"
$ws.Range("E11").Value = "  -9.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("E12").Value = "  -4.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000216"
$ws.Range("E13").Value = "  -4.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.45"
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").Value = "3.449.32"
$ws.Range("E15").Value = "  -3.71%  "
$ws.Range("D16").Value = "61.236.50"
$ws.Range("E16").Value = "  -3.64%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "2.968.36"
$ws.Range("E18").Value = "  -3.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.54"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.92"
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.95"
$ws.Range("E21").Value = "  -4.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.657"
$ws.Range("E22").Value = "  -6.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.90"
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.34"
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.90"
$ws.Range("E25").Value = "  -2.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.67"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.56"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.87"
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.26"
$ws.Range("E31").Value = "  -3.93%  "
$ws.Range("E32").Value = "  -3.95%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.26"
$ws.Range("E33").Value = "  -2.64%  "
$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.40"
$ws.Range("E34").Value = "  -1.39%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "54.50"
$ws.Range("E35").Value = "  -6.35%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.81"
$ws.Range("E36").Value = "  -3.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "444.42"
$ws.Range("E37").Value = "  -9.84%  "
$ws.Range("D38").Value = "3.088.66"
$ws.Range("E38").Value = "  -5.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0782"
$ws.Range("E39").Value = "  -2.38%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0374"
$ws.Range("E40").Value = "  -7.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("E41").Value = "  -2.71%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.99"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.28"
$ws.Range("E44").Value = "  -12.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.07"
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.238"
$ws.Range("E46").Value = "  -6.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.107"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.91"
$ws.Range("E48").Value = "  -6.02%  "
$ws.Range("B49").Value = "BitgetToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.30"
$ws.Range("E49").Value = "  +8.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.98"
$ws.Range("E50").Value = "  -7.69%  "
$ws.Range("D51").Value = "0.0₃0473"
$ws.Range("E51").Value = "  -11.19%  "
